# Added Week 15 simulations
# Update the cumulative target-depth counts on both the "OFF" and "DEF"
# sheets (row 2 = "H" totals) to reflect the newly simulated week.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 206
$wsOff.Range("C2").Value = 144
$wsOff.Range("D2").Value = 52
$wsOff.Range("E2").Value = 26

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 189
$wsDef.Range("C2").Value = 131
$wsDef.Range("D2").Value = 41
$wsDef.Range("E2").Value = 21
$wsDef.Range("F2").Value = 4
